$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.728.82'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.39%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.289.44'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.09%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '96.25'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.91%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '269.66'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.10%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.617'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.24%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.09%  '

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.47%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.47'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.22%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0935'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.62%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.97'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.21%  '

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.78%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.69'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.88%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.632.02'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.01%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.854'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.46%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.293.00'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.97%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.699.57'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.36%  '

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.13%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.18'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.83%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.10'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.99%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.49'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +10.67%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '232.64'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.81%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.11'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -4.73%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.71'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +8.39%  '

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.12%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.29'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.11%  '

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.04%  '

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.17%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.85'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.85%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '174.73'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.93%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.80'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.67%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0896'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.47%  '

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.87%  '

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.35%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.53'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +4.35%  '

$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.107'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.42%  '

$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0351'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.42%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.57'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +4.89%  '

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.96%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.31'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.03%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '12.37'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.96%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.34'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.28%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '64.22'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.49%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.73'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.11%  '

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.30%  '

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.13%  '

$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.20'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.55%  '

$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '97.40'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.68%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.51'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +11.25%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.431'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.63%  '
